$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 469, pushing the existing rows 469:492 down to 471:494.
$ws.Rows.Item(469).Insert()
$ws.Rows.Item(469).Insert()

# New weekly entry (Primera) for date 45267 (2023-12-07)
$ws.Range("A469").Value = 11
$ws.Range("B469").Value = "Vega Monumental Concepción"
$ws.Range("C469").Value = "Bíobío"
$ws.Range("D469").Value = 45267
$ws.Range("E469").Value = 8
$ws.Range("F469").Value = 100112009
$ws.Range("G469").Value = "Acelga"
$ws.Range("H469").Value = "Sin especificar"
$ws.Range("I469").Value = "Primera"
$ws.Range("J469").Value = 200
$ws.Range("K469").Value = 700
$ws.Range("L469").Value = 800
$ws.Range("M469").Value = 750
$ws.Range("N469").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O469").Value = "Región de Ñuble"
$ws.Range("P469").Value = 750
$ws.Range("Q469").Value = 1
$ws.Range("R469").Value = "Hortaliza"

# New weekly entry (Segunda) for date 45267 (2023-12-07)
$ws.Range("A470").Value = 11
$ws.Range("B470").Value = "Vega Monumental Concepción"
$ws.Range("C470").Value = "Bíobío"
$ws.Range("D470").Value = 45267
$ws.Range("E470").Value = 8
$ws.Range("F470").Value = 100112009
$ws.Range("G470").Value = "Acelga"
$ws.Range("H470").Value = "Sin especificar"
$ws.Range("I470").Value = "Segunda"
$ws.Range("J470").Value = 100
$ws.Range("K470").Value = 600
$ws.Range("L470").Value = 600
$ws.Range("M470").Value = 600
$ws.Range("N470").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O470").Value = "Región de Ñuble"
$ws.Range("P470").Value = 600
$ws.Range("Q470").Value = 1
$ws.Range("R470").Value = "Hortaliza"
